# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / *NQ / *HQ / LevePrice* / LeveProfit*) on each
# class sheet's leve-profit table with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 228.88235
$ws.Range("I33").Value = 213.23077
$ws.Range("K33").Value = 213.23077
$ws.Range("M33").Value = 15.76922999999999

$ws.Range("H64").Value = 3461.087
$ws.Range("I64").Value = 3406.25
$ws.Range("J64").Value = 3586.4285
$ws.Range("K64").Value = 3406.25
$ws.Range("L64").Value = 3586.4285
$ws.Range("M64").Value = -3158.25
$ws.Range("N64").Value = -4082.4285

$ws.Range("H67").Value = 3461.087
$ws.Range("I67").Value = 3406.25
$ws.Range("J67").Value = 3586.4285
$ws.Range("K67").Value = 3406.25
$ws.Range("L67").Value = 3586.4285
$ws.Range("M67").Value = -2548.25
$ws.Range("N67").Value = -5302.4285

$ws.Range("H74").Value = 6371.143
$ws.Range("I74").Value = 4959.6
$ws.Range("J74").Value = 9900
$ws.Range("K74").Value = 4959.6
$ws.Range("L74").Value = 9900
$ws.Range("N74").Value = -11772

$ws.Range("H77").Value = 6371.143
$ws.Range("I77").Value = 4959.6
$ws.Range("J77").Value = 9900
$ws.Range("K77").Value = 24798
$ws.Range("L77").Value = 49500
$ws.Range("N77").Value = -58860

$ws.Range("H100").Value = 15386350
$ws.Range("I100").Value = 15386350
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 15386350
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -15385809
$ws.Range("N100").ClearContents()

$ws.Range("H106").Value = 3862.9333
$ws.Range("I106").Value = 1914.4
$ws.Range("K106").Value = 1914.4
$ws.Range("M106").Value = -1283.4

$ws.Range("H112").Value = 26317180
$ws.Range("I112").Value = 250000460
$ws.Range("K112").Value = 750001380
$ws.Range("M112").Value = -750000272

$ws.Range("H113").Value = 3381.2727
$ws.Range("I113").Value = 1372.25
$ws.Range("J113").Value = 3827.7222
$ws.Range("K113").Value = 1372.25
$ws.Range("L113").Value = 3827.7222
$ws.Range("M113").Value = 1881.75
$ws.Range("N113").Value = -10335.7222

$ws.Range("H129").Value = 819.88
$ws.Range("J129").Value = 899.3488
$ws.Range("L129").Value = 2698.0464
$ws.Range("N129").Value = -12698.0464

$ws.Range("H132").Value = 218804.95
$ws.Range("I132").Value = 1305.4839
$ws.Range("J132").Value = 668303.9
$ws.Range("K132").Value = 3916.4517
$ws.Range("L132").Value = 2004911.7
$ws.Range("M132").Value = -1386.4517
$ws.Range("N132").Value = -2009971.7

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H137").Value = 1135596.6
$ws.Range("I137").Value = 1701885.1
$ws.Range("J137").Value = 3019.5715
$ws.Range("K137").Value = 5105655.300000001
$ws.Range("L137").Value = 9058.7145
$ws.Range("M137").Value = -5103105.300000001
$ws.Range("N137").Value = -14158.7145

$ws.Range("H138").Value = 7833.8
$ws.Range("I138").Value = 1584.4445
$ws.Range("J138").Value = 8451.868
$ws.Range("K138").Value = 4753.333500000001
$ws.Range("L138").Value = 25355.604
$ws.Range("M138").Value = 386.6664999999994
$ws.Range("N138").Value = -35635.604

$ws.Range("H141").Value = 35959.438
$ws.Range("I141").Value = 43404.04
$ws.Range("K141").Value = 130212.12
$ws.Range("M141").Value = -125032.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3820.8481
$ws.Range("I32").Value = 3183.7122
$ws.Range("K32").Value = 3183.7122
$ws.Range("M32").Value = -2896.7122

$ws.Range("H74").Value = 4539.759
$ws.Range("I74").Value = 5066.3687
$ws.Range("K74").Value = 5066.3687
$ws.Range("M74").Value = -4192.3687

$ws.Range("H77").Value = 4539.759
$ws.Range("I77").Value = 5066.3687
$ws.Range("K77").Value = 25331.8435
$ws.Range("M77").Value = -20963.8435

$ws.Range("H103").Value = 34288.75
$ws.Range("J103").Value = 34288.75
$ws.Range("L103").Value = 34288.75
$ws.Range("N103").Value = -36632.75

$ws.Range("H132").Value = 1585.597
$ws.Range("I132").Value = 1043.1395
$ws.Range("K132").Value = 3129.4185
$ws.Range("M132").Value = -599.4184999999998

$ws.Range("H133").Value = 36097.1
$ws.Range("J133").Value = 36097.1
$ws.Range("L133").Value = 36097.1
$ws.Range("N133").Value = -41157.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1634.2623
$ws.Range("I105").Value = 1635.5593
$ws.Range("K105").Value = 1635.5593
$ws.Range("M105").Value = 111.4407000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2920.8667
$ws.Range("I31").Value = 1090.421
$ws.Range("J31").Value = 6082.5454
$ws.Range("K31").Value = 1090.421
$ws.Range("L31").Value = 6082.5454
$ws.Range("M31").Value = -795.421
$ws.Range("N31").Value = -6672.5454

$ws.Range("H34").Value = 2920.8667
$ws.Range("I34").Value = 1090.421
$ws.Range("J34").Value = 6082.5454
$ws.Range("K34").Value = 1090.421
$ws.Range("L34").Value = 6082.5454
$ws.Range("M34").Value = -888.421
$ws.Range("N34").Value = -6486.5454

$ws.Range("H122").Value = 2129.1667
$ws.Range("I122").Value = 1026.8
$ws.Range("J122").Value = 3507.125
$ws.Range("K122").Value = 3080.4
$ws.Range("L122").Value = 10521.375
$ws.Range("M122").Value = -630.3999999999996
$ws.Range("N122").Value = -15421.375

$ws.Range("H132").Value = 2881.4
$ws.Range("I132").Value = 2453.1143
$ws.Range("K132").Value = 7359.342900000001
$ws.Range("M132").Value = -4829.342900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 554.9
$ws.Range("J113").Value = 565.4666999999999
$ws.Range("L113").Value = 1696.4001
$ws.Range("N113").Value = -6036.4001

$ws.Range("H131").Value = 5435579.5
$ws.Range("I131").Value = 55555964
$ws.Range("J131").Value = 839.0482
$ws.Range("K131").Value = 166667892
$ws.Range("L131").Value = 2517.1446
$ws.Range("M131").Value = -166662852
$ws.Range("N131").Value = -12597.1446

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2724.147
$ws.Range("I122").Value = 2075.2778
$ws.Range("J122").Value = 3454.125
$ws.Range("K122").Value = 6225.8334
$ws.Range("L122").Value = 10362.375
$ws.Range("M122").Value = -3775.8334
$ws.Range("N122").Value = -15262.375

$ws.Range("H126").Value = 2924.86
$ws.Range("I126").Value = 2820.9885
$ws.Range("J126").Value = 3620
$ws.Range("K126").Value = 8462.9655
$ws.Range("L126").Value = 10860
$ws.Range("M126").Value = -5992.9655
$ws.Range("N126").Value = -15800

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1552.5
$ws.Range("J46").Value = 1330.5333
$ws.Range("L46").Value = 1330.5333
$ws.Range("N46").Value = -1706.5333

$ws.Range("H55").Value = 254.46153
$ws.Range("I55").Value = 242.75
$ws.Range("J55").Value = 273.2
$ws.Range("K55").Value = 242.75
$ws.Range("L55").Value = 273.2
$ws.Range("M55").Value = -69.75
$ws.Range("N55").Value = -619.2

$ws.Range("H122").Value = 4997.778
$ws.Range("I122").Value = 2550
$ws.Range("J122").Value = 6956
$ws.Range("K122").Value = 7650
$ws.Range("L122").Value = 20868
$ws.Range("M122").Value = -5200
$ws.Range("N122").Value = -25768

$ws.Range("H136").Value = 3060.675
$ws.Range("I136").Value = 1435.9565
$ws.Range("J136").Value = 5258.8237
$ws.Range("K136").Value = 4307.8695
$ws.Range("L136").Value = 15776.4711
$ws.Range("M136").Value = -1757.8695
$ws.Range("N136").Value = -20876.4711

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 14700.25
$ws.Range("J74").Value = 14700.25
$ws.Range("L74").Value = 14700.25
$ws.Range("N74").Value = -16572.25

$ws.Range("H77").Value = 14700.25
$ws.Range("J77").Value = 14700.25
$ws.Range("L77").Value = 44100.75
$ws.Range("N77").Value = -53460.75

$ws.Range("H132").Value = 1710.7046
$ws.Range("I132").Value = 1185.3438
$ws.Range("J132").Value = 3111.6667
$ws.Range("K132").Value = 3556.0314
$ws.Range("L132").Value = 9335.000100000001
$ws.Range("M132").Value = -1026.0314
$ws.Range("N132").Value = -14395.0001

$ws.Range("H136").Value = 3219.9092
$ws.Range("I136").Value = 2094.2
$ws.Range("J136").Value = 6737.75
$ws.Range("K136").Value = 6282.599999999999
$ws.Range("L136").Value = 20213.25
$ws.Range("M136").Value = -3732.599999999999
$ws.Range("N136").Value = -25313.25
